# #2 Changing the return type to array into array
#
# Replace the single-column sample data (A1:A5 = "A".."E") with a small
# "people" table that has a header row (Name / Email / Ranking) and three
# data rows, each with a mailto: hyperlink on the Email column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Ranking"

# Data rows
$ws.Range("A2").Value = "josdem"
$ws.Range("B2").Value = "joseluis.delacruz@gmail.com"
$ws.Range("C2").Value = 5

$ws.Range("A3").Value = "eric"
$ws.Range("B3").Value = "erich@email.com"
$ws.Range("C3").Value = 5

$ws.Range("A4").Value = "martin"
$ws.Range("B4").Value = "martinv@email.com"
$ws.Range("C4").Value = 5

# Remove the old 5th row of sample data (A5 = "E") which no longer exists
# in the new, smaller table.
$ws.Range("A5").ClearContents()

# Turn the e-mail addresses into real hyperlinks (mailto:), which also
# applies the built-in "Hyperlink" cell style to each of them.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:joseluis.delacruz@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:erich@email.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:martinv@email.com")

# Resize the columns to fit their new contents.
$ws.Columns("A:C").AutoFit()

# Leave the selection where Excel would land after typing in the last
# data row (one row below, one column to the right of the table).
$ws.Range("D6").Select()
